$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated experiment results (new k for cross-validation).
# Row 2 = random_forest
$ws.Range("B2").Value = 4.1412002856448709
$ws.Range("C2").Value = 0.30546264058080946
$ws.Range("D2").Value = 3.3292254131338739
$ws.Range("E2").Value = 0.32806419305964213
$ws.Range("F2").Value = 0.57276888276131244
$ws.Range("G2").Value = 0.54354766190198689
$ws.Range("H2").Value = 0.67193580694035782
$ws.Range("I2").Value = 0.83880020743338546

# Row 3 = lsboost
$ws.Range("B3").Value = 4.485909705210414
$ws.Range("C3").Value = 0.33088904893361681
$ws.Range("D3").Value = 3.5058552321211351
$ws.Range("E3").Value = 0.38495274078378017
$ws.Range("F3").Value = 0.62044559856910919
$ws.Range("G3").Value = 0.57238521815574739
$ws.Range("H3").Value = 0.61504725921621983
$ws.Range("I3").Value = 0.78743246347649398

# Row 4 = neural_network
$ws.Range("B4").Value = 3.8997687768617166
$ws.Range("C4").Value = 0.28765420314590595
$ws.Range("D4").Value = 3.1808261499344175
$ws.Range("E4").Value = 0.29092703206125475
$ws.Range("F4").Value = 0.53937652160735983
$ws.Range("G4").Value = 0.51931918154081091
$ws.Range("H4").Value = 0.7090729679387453
$ws.Range("I4").Value = 0.86010400609229887

$wb.Save()
